$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 9 (NYKU, 517572) first so row indices above it stay stable
$ws.Rows.Item(9).Delete()

# Delete row 5 (FFAU, 142477)
$ws.Rows.Item(5).Delete()

# After deletions, rows shift up:
# old row6 (ONEU) -> new row5
# old row7 (GAOU) -> new row6
# old row8 (FDCU) -> new row7
# Update "Carga limpia" (col G) to "NO" for TCLU (row4), ONEU (row5), GAOU (row6), FDCU (row7)
$ws.Range("G4").Value = "NO"
$ws.Range("G5").Value = "NO"
$ws.Range("G6").Value = "NO"
$ws.Range("G7").Value = "NO"
